$wb = $excel.ActiveWorkbook

# --- Summary sheet: just update the selection ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("C6").Select()

# --- Transactions sheet: update values and selection ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
$wsTransactions.Range("A2").Value = 75
$wsTransactions.Range("A3").Value = 74
$wsTransactions.Range("D3").Select()
